# NC-3.2-3.3 - Rotina para retornar todos os crimes com filtro de datas
# + Plotar no app as ocorrências nas localizações onde ocorreram
#
# This script reproduces, via Excel COM automation, the changes described
# by the target diff:
#  1) Tasks 3.2 (row 19) and 3.3 (row 20) are marked "Concluído" (green),
#     matching the style already used by other finished tasks.
#  2) Task 4.6 (row 28, "cadastro da ocorrência de agressão") moves from
#     "Fazendo" (yellow) to "Concluído" (green).
#  3) A new task 4.7 ("cadastro da ocorrência outros") is appended to user
#     story 4, as a new row with no status set yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$green = 5296274   # RGB 92D050 (BGR-encoded OLE color) - used for "done" rows

# ---------------------------------------------------------------------
# 1) Task 3.2 (row 19) -> Concluído
# ---------------------------------------------------------------------
$rng = $ws.Range("C19:D19")
$rng.Interior.Color = $green
$rng.WrapText = $true

$status = $ws.Range("E19")
$status.Value2 = "Concluído"
$status.Interior.Color = $green
$status.WrapText = $false

# ---------------------------------------------------------------------
# 2) Task 3.3 (row 20) -> Concluido
# ---------------------------------------------------------------------
$rng = $ws.Range("C20:D20")
$rng.Interior.Color = $green
$rng.WrapText = $true

$status = $ws.Range("E20")
$status.Value2 = "Concluido"
$status.Interior.Color = $green
$status.WrapText = $false

# ---------------------------------------------------------------------
# 3) Insert a new row inside the US-4 group (rows 23:28) so the merged
#    A23/B23 cells grow to 23:29, then populate it as task 4.7.
# ---------------------------------------------------------------------
$ws.Rows.Item(28).Insert()

# The former row 28 (task 4.6, "Fazendo"/yellow) is now row 29 - turn it
# into the brand-new task 4.7 with a plain (no fill) style and no status.
$c29 = $ws.Range("C29")
$c29.Value2 = "Criar rotina no backend para cadastro da ocorrência outros"
$c29.Interior.ColorIndex = -4142
$c29.WrapText = $true

$d29 = $ws.Range("D29")
$d29.Value2 = "4.7"
$d29.Interior.ColorIndex = -4142
$d29.WrapText = $true

$e29 = $ws.Range("E29")
$e29.Value2 = ""
$e29.Interior.ColorIndex = -4142
$e29.WrapText = $false

# The freshly inserted blank row 28 inherited its neighbour's green style;
# fill it back in with the original task-4.6 text now marked Concluído.
$c28 = $ws.Range("C28")
$c28.Value2 = "Criar rotina no backend para cadastro da ocorrência de agressão"
$c28.Interior.Color = $green
$c28.WrapText = $true

$d28 = $ws.Range("D28")
$d28.Value2 = "4.6"
$d28.Interior.Color = $green
$d28.WrapText = $true

$e28 = $ws.Range("E28")
$e28.Value2 = "Concluído"
$e28.Interior.Color = $green
$e28.WrapText = $false

# ---------------------------------------------------------------------
# 4) Selection / scroll position, matching what Excel would leave behind
#    after the user finished editing around the newly-inserted row.
# ---------------------------------------------------------------------
$ws.Application.GoTo($ws.Range("B38:B41"), $true)
$ws.Range("B38:B41").Select()
$ws.Application.ActiveWindow.ScrollRow = 23
